$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" footer field text (16/05/2019 ->
#    19/10/2022) on the slide master and every slide layout's Date
#    Placeholder. ppPlaceholderDate = 16.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Type -eq 14) {
            $pt = $shp.PlaceholderFormat.Type
            if ($pt -eq 16) {
                $shp.TextFrame.TextRange.Text = "19/10/2022"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Re-embed / touch the Visio OLE drawing on the "File format changes:"
#    slide so the legacy VML shape id backing p:oleObj/@spid is reallocated
#    (_x0000_s3087 -> _x0000_s3088), matching a PowerPoint re-save of the
#    embedded object.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.Type -eq 7) {
            $of = $shp.OLEFormat
            if ($of -ne $null) {
                $of.Activate()
                $of.DoVerb(0)
            }
        }
    }
}
